# Apply updated cryptocurrency market data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: assign directly ---
$ws.Range('D2').Value = '42.963.54'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.302.87'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E6').Value = '  -5.01%  '
$ws.Range('E7').Value = '  -1.38%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -3.31%  '
$ws.Range('E10').Value = '  -4.58%  '
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('E12').Value = '  -4.58%  '
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('E14').Value = '  +8.02%  '
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '2.661.17'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').Value = '2.288.03'
$ws.Range('E17').Value = '  -3.64%  '
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').Value = '42.901.33'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0901'
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E21').Value = '  -2.46%  '
$ws.Range('E22').Value = '  -2.09%  '
$ws.Range('E23').Value = '  -1.68%  '
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -3.56%  '
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('E29').Value = '  +3.83%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('E31').Value = '  -2.39%  '
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  +6.25%  '
$ws.Range('E35').Value = '  -2.25%  '
$ws.Range('E36').Value = '  -1.05%  '
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('E39').Value = '  -3.35%  '
$ws.Range('E40').Value = '  -2.34%  '
$ws.Range('E41').Value = '  -4.48%  '
$ws.Range('E42').Value = '  -2.09%  '
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('D44').Value = '1.975.42'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('E45').Value = '  -1.91%  '
$ws.Range('E46').Value = '  -3.84%  '
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('E48').Value = '  -3.73%  '
$ws.Range('D49').Value = '2.528.85'
$ws.Range('E50').Value = '  -6.99%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('E51').Value = '  -6.39%  '

# --- Numeric-looking text values (e.g. "2.00", "301.01"): Excel would
# auto-convert a plain .Value assignment into a floating point number and
# lose the exact original text (trailing zeros / precision). To preserve
# them as literal text, stage each value in a helper cell that has been
# explicitly formatted as Text, then copy/paste-special just the value
# into the destination cell (this keeps the destination cell style untouched).
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"
$helper.Value = '301.01'
$helper.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$helper.Value = '97.05'
$helper.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$helper.Value = '33.68'
$helper.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$helper.Value = '0.0794'
$helper.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$helper.Value = '49.28'
$helper.Copy()
$ws.Range('D12').PasteSpecial(-4163)
$helper.Value = '16.85'
$helper.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$helper.Value = '0.808'
$helper.Copy()
$ws.Range('D18').PasteSpecial(-4163)
$helper.Value = '11.59'
$helper.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$helper.Value = '67.19'
$helper.Copy()
$ws.Range('D23').PasteSpecial(-4163)
$helper.Value = '236.31'
$helper.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$helper.Value = '2.00'
$helper.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$helper.Value = '2.46'
$helper.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$helper.Value = '24.82'
$helper.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$helper.Value = '166.49'
$helper.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$helper.Value = '33.85'
$helper.Copy()
$ws.Range('D31').PasteSpecial(-4163)
$helper.Value = '9.12'
$helper.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$helper.Value = '4.79'
$helper.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$helper.Value = '2.39'
$helper.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$helper.Value = '0.0696'
$helper.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$helper.Value = '2.35'
$helper.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$helper.Value = '17.77'
$helper.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$helper.Value = '9.81'
$helper.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$helper.Value = '52.78'
$helper.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$helper.Value = '4.58'
$helper.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = 0
